# "#5: fund, bonds, otherbonds, antique done"
#
# This batch of edits normalizes the "基金受益憑證" (fund) sheet into the
# same flattened record layout already used by the other sheets (stock,
# deposit, etc.) and removes the still-unprocessed "其他有價證券"
# (other securities / antiques) sheet entirely, since its rows hadn't
# been normalized yet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Drop the not-yet-normalized "其他有價證券" sheet. Excel automatically
#    re-targets the relationship ids of the remaining sheets (債務 keeps
#    rId6, 事業投資 keeps rId7) once this sheet is removed.
# ---------------------------------------------------------------------
$otherSecurities = $wb.Worksheets.Item("其他有價證券")
$otherSecurities.Delete()

# ---------------------------------------------------------------------
# 2. Rebuild "基金受益憑證" (fund) into the normalized layout:
#    - row 1 becomes a real header row (field names) instead of
#      duplicating the first record
#    - existing columns B:H keep their values, just shifted onto the
#      header-described columns
#    - new columns I:O carry the normalized metadata for every record
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("基金受益憑證")

# Capture the existing B:H values for data rows 2-9 before overwriting
# row 1 (row 1 currently duplicates row 2's data).
$lastRow = 9
$data = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $data[$r] = @(
        $ws.Cells.Item($r, 2).Value(),
        $ws.Cells.Item($r, 3).Value(),
        $ws.Cells.Item($r, 4).Value(),
        $ws.Cells.Item($r, 5).Value(),
        $ws.Cells.Item($r, 6).Value(),
        $ws.Cells.Item($r, 7).Value(),
        $ws.Cells.Item($r, 8).Value()
    )
}

# Clone the bold/bordered header style (already present on H1) across the
# new header cells, and the plain data style (already present on H2:H9)
# across the new data columns, before writing values into them.
$ws.Range("H1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122)
$ws.Range("H2:H9").Copy()
$ws.Range("I2:O9").PasteSpecial(-4122)

# Row 1: header labels.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# Rows 2-9: keep original B:H, fill in the normalized I:O metadata.
for ($r = 2; $r -le $lastRow; $r++) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
    $ws.Cells.Item($r, 8).Value = $vals[6]

    $ws.Cells.Item($r, 9).Value = "fund"
    $ws.Cells.Item($r, 10).Value = "normal"
    $ws.Cells.Item($r, 11).Value = "2012-04-30"
    $ws.Cells.Item($r, 12).Value = "徐少萍"
    $ws.Cells.Item($r, 13).Value = 726
    $ws.Cells.Item($r, 14).Value = "tmpf37d1"
    $ws.Cells.Item($r, 15).Value = $ws.Cells.Item($r, 1).Value()
}

Write-Host "done"
